$wb = $excel.ActiveWorkbook

# --- 1. Update "Nädal 10" (sheet 10): a couple of logged time entries changed ---
$ws10 = $wb.Worksheets.Item(10)
$ws10.Range("D19").Value = 0.0027777777777777779
$ws10.Range("F19").Value = 84

# Update the comment text for week 10 (string shared with other "tehtud" notes)
$ws10.Range("H19").Value = "p. 56, 57 tehtud"

# --- 2. Duplicate "Nädal 10" to create the new "Nädal 11" sheet ---
$ws10.Copy([System.Reflection.Missing]::Value, $ws10)
$ws11 = $wb.Worksheets.Item(11)
$ws11.Name = "Nädal 11"

# New week's blank template: clear the logged entries (B:J, rows 7-19) while
# keeping row numbering (col A) and all formatting/styles intact.
$ws11.Range("B7:J19").ClearContents()

# New week's date-range header
$ws11.Range("G4").Value = "07.04.2020 - 13.04.2020"

# New sheet's selection (per template default) and make sure it is not the
# tab shown as active/selected — "Nädal 10" stays the active tab.
$ws11.Range("H15").Select()

# --- 3. Restore "Nädal 10" as the active/selected sheet + its own selection ---
$ws10.Activate()
$ws10.Range("H16").Select()
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(1, 1)
